$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) contains numeric-looking text (e.g. "273.78", "1.800").
# Force those specific cells to Text format before assigning so Excel keeps
# them as literal strings instead of coercing them into numbers, matching the
# original inlineStr cell contents. Other cells/formatting are left untouched.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = '19.974.27'
$ws.Range("E2").Value = '  -8.06%  '
$ws.Range("D3").Value = '1.417.39'
$ws.Range("E3").Value = '  -7.87%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("D6").Value = '273.78'
$ws.Range("E6").Value = '  -5.43%  '
$ws.Range("D7").Value = '0.3693'
$ws.Range("E7").Value = '  -5.10%  '
$ws.Range("D8").Value = '0.3078'
$ws.Range("E8").Value = '  -3.37%  '
$ws.Range("D9").Value = '39.57'
$ws.Range("E9").Value = '  -7.77%  '
$ws.Range("D10").Value = '1.003'
$ws.Range("E10").Value = '  -5.05%  '
$ws.Range("D11").Value = '0.06583'
$ws.Range("E11").Value = '  -8.57%  '
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").Value = '5.445'
$ws.Range("E13").Value = '  -3.38%  '
$ws.Range("D14").Value = '17.06'
$ws.Range("E14").Value = '  -8.15%  '
$ws.Range("D15").Value = '6.166'
$ws.Range("E15").Value = '  -6.70%  '
$ws.Range("D16").Value = '1.426.52'
$ws.Range("E16").Value = '  -7.58%  '
$ws.Range("E17").Value = '  -9.36%  '
$ws.Range("D18").Value = '0.05750'
$ws.Range("E18").Value = '  -12.67%  '
$ws.Range("D19").Value = '74.21'
$ws.Range("E19").Value = '  -10.77%  '
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").Value = '5.624'
$ws.Range("E21").Value = '  -8.42%  '
$ws.Range("D22").Value = '14.49'
$ws.Range("E22").Value = '  -5.74%  '
$ws.Range("D23").Value = '10.96'
$ws.Range("E23").Value = '  +1.01%  '
$ws.Range("D24").Value = '2.335'
$ws.Range("E24").Value = '  -2.89%  '
$ws.Range("D25").Value = '19.991.97'
$ws.Range("E25").Value = '  -7.98%  '
$ws.Range("D26").Value = '2.278'
$ws.Range("E26").Value = '  -3.82%  '
$ws.Range("D27").Value = '139.57'
$ws.Range("E27").Value = '  -4.50%  '
$ws.Range("D28").Value = '16.96'
$ws.Range("E28").Value = '  -7.61%  '
$ws.Range("D29").Value = '1.582.16'
$ws.Range("E29").Value = '  -7.77%  '
$ws.Range("D30").Value = '109.20'
$ws.Range("E30").Value = '  -7.15%  '
$ws.Range("D31").Value = '3.887'
$ws.Range("E31").Value = '  -19.59%  '
$ws.Range("D32").Value = '5.398'
$ws.Range("E32").Value = '  -8.60%  '
$ws.Range("D33").Value = '0.8555'
$ws.Range("E33").Value = '  -11.77%  '
$ws.Range("D34").Value = '0.07725'
$ws.Range("E34").Value = '  -5.59%  '
$ws.Range("D35").Value = '8.438'
$ws.Range("E35").Value = '  -4.21%  '
$ws.Range("D36").Value = '0.05744'
$ws.Range("E36").Value = '  -5.63%  '
$ws.Range("D37").Value = '4.779'
$ws.Range("E37").Value = '  -6.85%  '
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("D39").Value = '10.66'
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = '0.1928'
$ws.Range("E40").Value = '  -5.33%  '
$ws.Range("D41").Value = '0.02034'
$ws.Range("E41").Value = '  -7.51%  '
$ws.Range("E42").Value = '  -10.12%  '
$ws.Range("D43").Value = '1.279'
$ws.Range("E43").Value = '  -14.05%  '
$ws.Range("D44").Value = '0.5308'
$ws.Range("E44").Value = '  -7.53%  '
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").Value = '3.532'
$ws.Range("E45").Value = '  -5.55%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '12.25'
$ws.Range("E46").Value = '  -6.44%  '
$ws.Range("E47").Value = '  -6.68%  '
$ws.Range("D48").Value = '1.800'
$ws.Range("E48").Value = '  -3.53%  '
$ws.Range("D49").Value = '109.46'
$ws.Range("E49").Value = '  -6.76%  '
$ws.Range("E50").Value = '  -10.06%  '
$ws.Range("E51").Value = '  +0.08%  '
